$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.979.04"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.642.45"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +2.38%  "
$ws.Range("D12").Value = "1.875.47"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.640.77"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.572"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "27.978.99"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "1.409.46"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +7.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "1.784.26"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
